# Gantt chart updates, Flowchart added, Moon and Spacestation copies
#
# Applies the task-table updates to the single "Gantt Chart" worksheet:
#   - Row 16 (2.1 Setup):            start date, % complete, fills "X" across AM:AQ
#   - Row 18 (2.3 Sounds):           owner "JC", start date, % complete, fills "X" across X:AE
#   - Row 19 (Lighting):             owner "JC", start date, % complete, fills "X" across AF:AL
#   - Row 21 (2.5 Map/Level Layout): owner "JC", start date, % complete, fills "X" across S:W
#   - Row 43 (4.1 new "Sliding Doors" task): owner, start/due dates, % complete, fills "X" across S:T

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# ---- Row 16 : WBS 2.1 "Setup" --------------------------------------------
$ws.Range("E16").Value = 44851
$ws.Range("H16").Value = 0.75
$ws.Range("H16").NumberFormat = "0%"
$ws.Range("H16").HorizontalAlignment = -4108
$ws.Range("H16").Borders.LineStyle = -4142
$ws.Range("AM16:AQ16").Value = "X"

# ---- Row 18 : WBS 2.3 "Sounds" -------------------------------------------
$ws.Range("D18").Value = "JC"
$ws.Range("E18").Value = 44830
$ws.Range("H18").Value = 0.5
$ws.Range("H18").NumberFormat = "0%"
$ws.Range("H18").HorizontalAlignment = -4108
$ws.Range("X18:AE18").Value = "X"

# ---- Row 19 : "Lighting" --------------------------------------------------
$ws.Range("D19").Value = "JC"
$ws.Range("E19").Value = 44840
$ws.Range("H19").Value = 0.25
$ws.Range("H19").NumberFormat = "0%"
$ws.Range("H19").HorizontalAlignment = -4108
$ws.Range("AF19:AL19").Value = "X"

# ---- Row 21 : WBS 2.5 "Map/Level Layout" ---------------------------------
$ws.Range("D21").Value = "JC"
$ws.Range("E21").Value = 44823
$ws.Range("H21").Value = 0.25
$ws.Range("H21").NumberFormat = "0%"
$ws.Range("H21").HorizontalAlignment = -4108
$ws.Range("S21:W21").Value = "X"

# ---- Row 43 : WBS 4.1 new task "Sliding Doors" ---------------------------
$ws.Range("C43").Value = "Sliding Doors"
$ws.Range("D43").Value = "JC"
$ws.Range("E43").Value = 44823
$ws.Range("F43").Value = 44825
$ws.Range("H43").Value = 1
$ws.Range("S43:T43").Value = "X"

# ---- Cursor position left where the author's last save left it ----------
$ws.Range("G24").Select()
